# Generate Report for Handoff
#
# Updates the localization-status report for the rows that were re-handed-off
# (1b3907e8-85b5-4b23-a8ad-c8a9c25ba682.md), bumping the recorded handoff/
# generation timestamps and stamping the new "ht" (handoff type) priority
# value on the per-language sheets.

$wb = $excel.ActiveWorkbook

$rows = 8, 9, 10, 12, 13, 14

# --- "Overview" sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-06 00:24:22"
}

# --- "zh-cn" sheet: column H = "Latest Handoff Datetime", column E = "Priority" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-09-06 00:24:16"
    $wsZhCn.Range("E$r").Value = "ht"
}

# --- "de-de" sheet: column E = "Priority" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
}
